# REV 1 - WIP
# Add a "Total" column (G) = Quantity (F) * 5, and update selection / column width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column G
$ws.Range("G1").Value = "Total"

# Fill G2:G46 with Quantity * 5 (computed from column F, written as static values)
for ($row = 2; $row -le 46; $row++) {
    $qty = $ws.Cells.Item($row, 6).Value2
    $ws.Cells.Item($row, 7).Value = $qty * 5
}

# Set column D width to match the new layout (~17.66 chars wide)
$ws.Columns.Item(4).ColumnWidth = 16.8333333333333

# Update the active selection to L8
$ws.Range("L8").Select()
